$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 538
$ws.Range("I31").Value = 313.5
$ws.Range("J31").Value = 874.75
$ws.Range("K31").Value = 940.5
$ws.Range("L31").Value = 2624.25
$ws.Range("M31").Value = -710.5
$ws.Range("N31").Value = -3084.25
$ws.Range("H38").Value = 1298
$ws.Range("I38").Value = 108
$ws.Range("J38").Value = 4868
$ws.Range("K38").Value = 324
$ws.Range("L38").Value = 14604
$ws.Range("M38").Value = 48
$ws.Range("N38").Value = -15348
$ws.Range("H40").Value = 1103.303
$ws.Range("I40").Value = 1045.8636
$ws.Range("J40").Value = 1218.1818
$ws.Range("K40").Value = 1045.8636
$ws.Range("L40").Value = 1218.1818
$ws.Range("M40").Value = -870.8635999999999
$ws.Range("N40").Value = -1568.1818
$ws.Range("H46").Value = 2000
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H53").Value = 41667510
$ws.Range("I53").Value = 66667664
$ws.Range("K53").Value = 66667664
$ws.Range("M53").Value = -66667027
$ws.Range("H54").Value = 1000
$ws.Range("I54").Value = 1000
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 1000
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -514
$ws.Range("N54").ClearContents()
$ws.Range("H55").Value = 990.2
$ws.Range("I55").Value = 51
$ws.Range("K55").Value = 51
$ws.Range("M55").Value = 163
$ws.Range("H59").Value = 1791.8
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 1791.8
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 5375.4
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -6489.4
$ws.Range("H60").Value = 2000
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H116").Value = 7059.6
$ws.Range("I116").Value = 9798.916999999999
$ws.Range("J116").Value = 2950.625
$ws.Range("K116").Value = 9798.916999999999
$ws.Range("L116").Value = 2950.625
$ws.Range("M116").Value = -6356.916999999999
$ws.Range("N116").Value = -9834.625
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H60").Value = 4112.75
$ws.Range("I60").Value = 3183.6667
$ws.Range("J60").Value = 6900
$ws.Range("K60").Value = 3183.6667
$ws.Range("L60").Value = 6900
$ws.Range("M60").Value = -2450.6667
$ws.Range("N60").Value = -8366
$ws.Range("H61").Value = 2501273
$ws.Range("I61").Value = 3846935.2
$ws.Range("J61").Value = 2185.7144
$ws.Range("K61").Value = 3846935.2
$ws.Range("L61").Value = 2185.7144
$ws.Range("M61").Value = -3846723.2
$ws.Range("N61").Value = -2609.7144
$ws.Range("H74").Value = 638.8611
$ws.Range("I74").Value = 531.8570999999999
$ws.Range("J74").Value = 1013.375
$ws.Range("K74").Value = 531.8570999999999
$ws.Range("L74").Value = 1013.375
$ws.Range("M74").Value = 342.1429000000001
$ws.Range("N74").Value = -2761.375
$ws.Range("H77").Value = 638.8611
$ws.Range("I77").Value = 531.8570999999999
$ws.Range("J77").Value = 1013.375
$ws.Range("K77").Value = 2659.2855
$ws.Range("L77").Value = 5066.875
$ws.Range("M77").Value = 1708.7145
$ws.Range("N77").Value = -13802.875
$ws.Range("H136").Value = 2501273
$ws.Range("I136").Value = 3846935.2
$ws.Range("J136").Value = 2185.7144
$ws.Range("K136").Value = 11540805.6
$ws.Range("L136").Value = 6557.1432
$ws.Range("M136").Value = -11538255.6
$ws.Range("N136").Value = -11657.1432
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H134").Value = 6661.2974
$ws.Range("I134").Value = 2285.8125
$ws.Range("J134").Value = 34664.4
$ws.Range("K134").Value = 6857.4375
$ws.Range("L134").Value = 103993.2
$ws.Range("M134").Value = -4322.4375
$ws.Range("N134").Value = -109063.2
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2128.8918
$ws.Range("I31").Value = 1882.3
$ws.Range("J31").Value = 3185.7144
$ws.Range("K31").Value = 1882.3
$ws.Range("L31").Value = 3185.7144
$ws.Range("M31").Value = -1587.3
$ws.Range("N31").Value = -3775.7144
$ws.Range("H34").Value = 2128.8918
$ws.Range("I34").Value = 1882.3
$ws.Range("J34").Value = 3185.7144
$ws.Range("K34").Value = 1882.3
$ws.Range("L34").Value = 3185.7144
$ws.Range("M34").Value = -1680.3
$ws.Range("N34").Value = -3589.7144
$ws.Range("H51").Value = 9597.4
$ws.Range("J51").Value = 9474.25
$ws.Range("L51").Value = 9474.25
$ws.Range("N51").Value = -10946.25
$ws.Range("H59").Value = 13638.25
$ws.Range("J59").Value = 13638.25
$ws.Range("L59").Value = 13638.25
$ws.Range("N59").Value = -15928.25
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H61").Value = 9597.4
$ws.Range("J61").Value = 9474.25
$ws.Range("L61").Value = 9474.25
$ws.Range("N61").Value = -10170.25
$ws.Range("H64").Value = 9950
$ws.Range("J64").Value = 9950
$ws.Range("L64").Value = 9950
$ws.Range("N64").Value = -10446
$ws.Range("H67").Value = 9950
$ws.Range("J67").Value = 9950
$ws.Range("L67").Value = 9950
$ws.Range("N67").Value = -11666
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 3458.3333
$ws.Range("I54").Value = 1000
$ws.Range("J54").Value = 4687.5
$ws.Range("K54").Value = 3000
$ws.Range("L54").Value = 14062.5
$ws.Range("M54").Value = -2441
$ws.Range("N54").Value = -15180.5
$ws.Range("H55").Value = 3154.9
$ws.Range("I55").Value = 675
$ws.Range("J55").Value = 3774.875
$ws.Range("K55").Value = 2025
$ws.Range("L55").Value = 11324.625
$ws.Range("M55").Value = -1848
$ws.Range("N55").Value = -11678.625
$ws.Range("H59").Value = 1000000000
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("H60").Value = 77777
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 77777
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 233331
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -233833
$ws.Range("H61").Value = 676.8570999999999
$ws.Range("I61").Value = 186.6
$ws.Range("J61").Value = 1902.5
$ws.Range("K61").Value = 559.8
$ws.Range("L61").Value = 5707.5
$ws.Range("M61").Value = -344.8
$ws.Range("N61").Value = -6137.5
$ws.Range("H113").Value = 1065.3158
$ws.Range("I113").Value = 990.4
$ws.Range("J113").Value = 1076.6666
$ws.Range("K113").Value = 2971.2
$ws.Range("L113").Value = 3229.9998
$ws.Range("M113").Value = -801.1999999999998
$ws.Range("N113").Value = -7569.9998
$ws.Range("H122").Value = 652.2
$ws.Range("I122").Value = 431.14285
$ws.Range("J122").Value = 1168
$ws.Range("K122").Value = 3880.28565
$ws.Range("L122").Value = 10512
$ws.Range("M122").Value = -1430.28565
$ws.Range("N122").Value = -15412
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7592.1875
$ws.Range("I46").Value = 1497.875
$ws.Range("J46").Value = 13686.5
$ws.Range("K46").Value = 1497.875
$ws.Range("L46").Value = 13686.5
$ws.Range("M46").Value = -1309.875
$ws.Range("N46").Value = -14062.5
$ws.Range("H60").Value = 6000
$ws.Range("J60").Value = 6000
$ws.Range("L60").Value = 6000
$ws.Range("N60").Value = -7018
$ws.Range("H132").Value = 3298
$ws.Range("I132").Value = 4899.7383
$ws.Range("J132").Value = 1319.3823
$ws.Range("K132").Value = 14699.2149
$ws.Range("L132").Value = 3958.1469
$ws.Range("M132").Value = -12169.2149
$ws.Range("N132").Value = -9018.1469
$ws.Range("H136").Value = 11376.706
$ws.Range("I136").Value = 10954.154
$ws.Range("J136").Value = 12750
$ws.Range("K136").Value = 32862.462
$ws.Range("L136").Value = 38250
$ws.Range("M136").Value = -30312.462
$ws.Range("N136").Value = -43350
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 10000
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H132").Value = 24196990
$ws.Range("I132").Value = 33785170
$ws.Range("J132").Value = 2024335.4
$ws.Range("K132").Value = 101355510
$ws.Range("L132").Value = 6073006.199999999
$ws.Range("M132").Value = -101352980
$ws.Range("N132").Value = -6078066.199999999
$ws.Range("H136").Value = 5397.864
$ws.Range("I136").Value = 6890.5
$ws.Range("J136").Value = 1417.5
$ws.Range("K136").Value = 20671.5
$ws.Range("L136").Value = 4252.5
$ws.Range("M136").Value = -18121.5
$ws.Range("N136").Value = -9352.5
